$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Information")
$ws.Range("B1").Value = 11
